$d = $word.ActiveDocument

$replacements = @(
    @{old = "90×22="; new = "19×94="},
    @{old = "43×48="; new = "12×69="},
    @{old = "40×60="; new = "70×36="},
    @{old = "63×52="; new = "23×96="},
    @{old = "46×98="; new = "64×59="},
    @{old = "33×37="; new = "88×94="},
    @{old = "94×64="; new = "78×88="},
    @{old = "72×32="; new = "19×57="},
    @{old = "41×57="; new = "91×44="},
    @{old = "44×53="; new = "65×90="},
    @{old = "77×89="; new = "33×87="},
    @{old = "94×80="; new = "53×30="},
    @{old = "29×29="; new = "89×23="},
    @{old = "67×19="; new = "77×90="},
    @{old = "69×21="; new = "20×37="},
    @{old = "46×21="; new = "18×99="},
    @{old = "53×75="; new = "83×43="},
    @{old = "25×12="; new = "21×69="},
    @{old = "47×69="; new = "64×89="},
    @{old = "99×48="; new = "99×80="},
    @{old = "94×47="; new = "52×40="},
    @{old = "89×48="; new = "60×86="},
    @{old = "21×56="; new = "49×50="},
    @{old = "22×14="; new = "21×63="},
    @{old = "65×16="; new = "61×93="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
